$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells that receive numeric-looking strings
# (e.g. "1.009", "29.736.37") stay text, matching the source inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($cell in $priceCells) {
  $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.736.37"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "2.096.93"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "343.86"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("D8").Value = "0.4384"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("D9").Value = "52.69"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "0.09279"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").Value = "1.165"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "24.91"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "2.101.85"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "8.265"
$ws.Range("D15").Value = "6.756"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "99.59"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "0.00001156"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "20.80"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "0.06648"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").Value = "29.760.30"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "2.350.50"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").Value = "2.511"
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("D29").Value = "161.33"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "133.07"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "1.139"
$ws.Range("E31").Value = "  -6.76%  "
$ws.Range("D32").Value = "0.1048"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D34").Value = "6.167"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").Value = "3.938"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "6.273"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "10.20"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").Value = "0.02576"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "0.06719"
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").Value = "12.47"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").Value = "0.6898"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").Value = "0.2223"
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("D43").Value = "1.318"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").Value = "0.6748"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "14.35"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "2.320"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "0.00000000358"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("D48").Value = "3.610"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").Value = "1.220"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "82.03"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("E51").Value = "  -2.10%  "
